$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Station")

# Insert a new column before column B, shifting id/name/machines/skills/kwargs right by one
$ws.Range("B1").EntireColumn.Insert()

# Copy the header style from the (now shifted) C1 cell to the new B1 cell
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122) | Out-Null

# Set the new header value
$ws.Range("B1").Value = "env"
